$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching formatting of existing header cells (bold, thin box border, centered)
$h1 = $ws.Range("H1")
$h1.Value = "Save"
$h1.Font.Bold = $true
$h1.Borders.Item("EdgeTop").LineStyle = "Continuous"
$h1.Borders.Item("EdgeTop").Weight = "Thin"
$h1.Borders.Item("EdgeBottom").LineStyle = "Continuous"
$h1.Borders.Item("EdgeBottom").Weight = "Thin"
$h1.Borders.Item("EdgeLeft").LineStyle = "Continuous"
$h1.Borders.Item("EdgeLeft").Weight = "Thin"
$h1.Borders.Item("EdgeRight").LineStyle = "Continuous"
$h1.Borders.Item("EdgeRight").Weight = "Thin"
$h1.HorizontalAlignment = "Center"
$h1.VerticalAlignment = "Top"

# Fill in the new "Save" column values for rows 2-7
$saveValues = @(0, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
